$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp text in F1
$ws.Range("F1").Value = "Last status check on: 14.02.2022 01:30"

# Fix row 5 (Makro): D5 and E5 were stored as text, convert to proper numeric values
$ws.Range("D5").Value = 0.6
$ws.Range("E5").Value = 44606.05216435185
$ws.Range("E5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
